# Refresh the cryptocurrency Price (column D) and Volume(1h) % change (column E)
# values for rows 2-51 with the latest pulled data (GitHub Actions data refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the "@" text storage the
# original cells use (inline/shared strings) even when the text looks like
# a number (e.g. "231.04") or contains only digits/dots, and restore the
# cell's original (default) style afterwards so no new styling is introduced.
function Set-TextValue($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "37.341.85"
Set-TextValue "E2" "  -1.33%  "
Set-TextValue "D3" "2.050.53"
Set-TextValue "E3" "  -1.56%  "
Set-TextValue "E4" "  -0.10%  "
Set-TextValue "D5" "231.04"
Set-TextValue "E5" "  -0.98%  "
Set-TextValue "D6" "0.622"
Set-TextValue "E6" "  -0.56%  "
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "57.01"
Set-TextValue "E8" "  -3.77%  "
Set-TextValue "E9" "  -2.90%  "
Set-TextValue "D10" "0.0770"
Set-TextValue "E10" "  -2.45%  "
Set-TextValue "E11" "  +1.28%  "
Set-TextValue "D12" "2.350.33"
Set-TextValue "E12" "  -1.67%  "
Set-TextValue "D13" "14.62"
Set-TextValue "E13" "  -1.02%  "
Set-TextValue "D14" "20.63"
Set-TextValue "E14" "  -2.70%  "
Set-TextValue "D15" "0.757"
Set-TextValue "E15" "  -2.53%  "
Set-TextValue "D16" "5.27"
Set-TextValue "E16" "  -1.78%  "
Set-TextValue "D17" "2.031.10"
Set-TextValue "E17" "  -1.91%  "
Set-TextValue "D18" "37.295.22"
Set-TextValue "E18" "  -1.23%  "
Set-TextValue "D19" "6.08"
Set-TextValue "E19" "  -1.42%  "
Set-TextValue "D20" "69.65"
Set-TextValue "E21" "  -3.23%  "
Set-TextValue "D22" "226.31"
Set-TextValue "E22" "  -0.85%  "
Set-TextValue "E23" "  +0.06%  "
Set-TextValue "E24" "  +0.00%  "
Set-TextValue "E25" "  -3.65%  "
Set-TextValue "D26" "9.82"
Set-TextValue "E26" "  +7.08%  "
Set-TextValue "D27" "169.99"
Set-TextValue "E27" "  -0.83%  "
Set-TextValue "E28" "  -5.96%  "
Set-TextValue "D29" "19.17"
Set-TextValue "E29" "  -1.64%  "
Set-TextValue "E30" "  -5.61%  "
Set-TextValue "E31" "  -0.12%  "
Set-TextValue "E32" "  -4.38%  "
Set-TextValue "E33" "  -1.70%  "
Set-TextValue "E34" "  -3.91%  "
Set-TextValue "E35" "  -1.61%  "
Set-TextValue "E36" "  +0.11%  "
Set-TextValue "D37" "3.28"
Set-TextValue "E37" "  -4.72%  "
Set-TextValue "E38" "  -0.08%  "
Set-TextValue "E39" "  -1.90%  "
Set-TextValue "E40" "  +3.15%  "
Set-TextValue "D41" "98.00"
Set-TextValue "E41" "  -1.09%  "
Set-TextValue "E42" "  -3.25%  "
Set-TextValue "E43" "  +0.08%  "
Set-TextValue "D44" "1.477.46"
Set-TextValue "E44" "  +2.26%  "
Set-TextValue "E45" "  +1.84%  "
Set-TextValue "D46" "16.59"
Set-TextValue "E47" "  -4.68%  "
Set-TextValue "E48" "  -3.13%  "
Set-TextValue "E49" "  -2.17%  "
Set-TextValue "D50" "2.93"
Set-TextValue "D51" "2.239.07"
Set-TextValue "E51" "  -1.59%  "
